# Update countries & provincias Spain
# Applies:
#  - Reorders "Dominica"/"Fiyi" and "Islas Virgenes Britanicas"/"Papua Nueva
#    Guinea" rows (the underlying per-country stats stay attached to the
#    country they belong to; only the label shown on each row swaps, mirroring
#    the upstream sharedStrings reorder).
#  - Refreshes the "Datos actualizados" timestamp cell.
#  - Refreshes the numeric covid stats for the countries whose figures moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Junio de 2020 a las 16:53"

# --- Country label swaps ---
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("A215").Value = "Islas Virgenes Britanicas"

# --- Updated country statistics ---
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2360593
$ws.Range("C4").Value = 3936
$ws.Range("E4").Value = 1257945
$ws.Range("G4").Value = 34
$ws.Range("H4").Value = 122281

# Row 8 - Reino Unido
$ws.Range("B8").Value = 305289
$ws.Range("C8").Value = 958
$ws.Range("G8").Value = 15
$ws.Range("H8").Value = 42647

# Row 14 - Alemania
$ws.Range("B14").Value = 191718
$ws.Range("C14").Value = 143
$ws.Range("E14").Value = 7854
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 8964

# Row 70 - Argelia
$ws.Range("E70").Value = 359
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 248

# Row 76 - Uzbekistan
$ws.Range("B76").Value = 6401
$ws.Range("C76").Value = 86
$ws.Range("D76").Value = 4450
$ws.Range("E76").Value = 1932

# Row 84 - Kenia
$ws.Range("B84").Value = 4797
$ws.Range("C84").Value = 59
$ws.Range("D84").Value = 1680
$ws.Range("E84").Value = 2992
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 125

# Row 102 - Cuba
$ws.Range("B102").Value = 2315
$ws.Range("C102").Value = 3
$ws.Range("D102").Value = 2113
$ws.Range("E102").Value = 117

# Row 137 - Principado de Andorra
$ws.Range("D137").Value = 796
$ws.Range("E137").Value = 7

# Row 158 - Mauricio
$ws.Range("B158").Value = 340
$ws.Range("C158").Value = 3
$ws.Range("E158").Value = 4

# Row 214 (now Papua Nueva Guinea)
$ws.Range("B214").Value = 9
$ws.Range("C214").Value = 1
$ws.Range("D214").Value = 8
$ws.Range("E214").Value = 1
$ws.Range("H214").Value = 0

# Row 215 (now Islas Virgenes Britanicas)
$ws.Range("D215").Value = 7
$ws.Range("H215").Value = 1
